$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-11-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-17 Monday", 2)

$table = $d.Tables.Item(1)

# Row 1
$table.Cell(1, 1).Range.Text = "522×8=4176"
$table.Cell(1, 2).Range.Text = "536×3=1608"
$table.Cell(1, 3).Range.Text = "837×3=2511"
$table.Cell(1, 4).Range.Text = "187×5=935"
$table.Cell(1, 5).Range.Text = "347×4=1388"

# Row 5
$table.Cell(5, 1).Range.Text = "803×9=7227"
$table.Cell(5, 2).Range.Text = "217×9=1953"
$table.Cell(5, 3).Range.Text = "777×6=4662"
$table.Cell(5, 4).Range.Text = "411×4=1644"
$table.Cell(5, 5).Range.Text = "475×9=4275"

# Row 10
$table.Cell(10, 1).Range.Text = "877×3=2631"
$table.Cell(10, 2).Range.Text = "148×6=888"
$table.Cell(10, 3).Range.Text = "561×9=5049"
$table.Cell(10, 4).Range.Text = "768×2=1536"
$table.Cell(10, 5).Range.Text = "657×9=5913"

# Row 15
$table.Cell(15, 1).Range.Text = "123×8=984"
$table.Cell(15, 2).Range.Text = "187×8=1496"
$table.Cell(15, 3).Range.Text = "433×9=3897"
$table.Cell(15, 4).Range.Text = "611×6=3666"
$table.Cell(15, 5).Range.Text = "438×8=3504"

# Row 20
$table.Cell(20, 1).Range.Text = "369×9=3321"
$table.Cell(20, 2).Range.Text = "276×9=2484"
$table.Cell(20, 3).Range.Text = "589×7=4123"
$table.Cell(20, 4).Range.Text = "860×6=5160"
$table.Cell(20, 5).Range.Text = "729×8=5832"
